$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, $text)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "54.501.41"
Set-TextValue $ws.Range("E2") "  -6.51%  "

Set-TextValue $ws.Range("D3") "2.444.50"
Set-TextValue $ws.Range("E3") "  -9.29%  "

Set-TextValue $ws.Range("E4") "  +0.00%  "

Set-TextValue $ws.Range("D5") "468.72"
Set-TextValue $ws.Range("E5") "  -6.06%  "

Set-TextValue $ws.Range("D6") "130.82"
Set-TextValue $ws.Range("E6") "  -5.63%  "

Set-TextValue $ws.Range("E7") "  +0.15%  "

Set-TextValue $ws.Range("D9") "2.440.62"
Set-TextValue $ws.Range("E9") "  -9.80%  "

Set-TextValue $ws.Range("E10") "  -8.80%  "

Set-TextValue $ws.Range("E11") "  -12.04%  "

Set-TextValue $ws.Range("E12") "  -8.77%  "

Set-TextValue $ws.Range("E13") "  -3.70%  "

Set-TextValue $ws.Range("D14") "2.870.78"
Set-TextValue $ws.Range("E14") "  -9.45%  "

Set-TextValue $ws.Range("D15") "54.571.42"
Set-TextValue $ws.Range("E15") "  -6.53%  "

Set-TextValue $ws.Range("E16") "  +0.14%  "

Set-TextValue $ws.Range("D17") "19.63"
Set-TextValue $ws.Range("E17") "  -8.00%  "

Set-TextValue $ws.Range("D18") "2.442.18"
Set-TextValue $ws.Range("E18") "  -9.66%  "

Set-TextValue $ws.Range("E19") "  -10.11%  "

Set-TextValue $ws.Range("D20") "311.95"
Set-TextValue $ws.Range("E20") "  -6.32%  "

Set-TextValue $ws.Range("D21") "9.56"
Set-TextValue $ws.Range("E21") "  -12.38%  "

Set-TextValue $ws.Range("E22") "  +0.27%  "

Set-TextValue $ws.Range("D23") "5.67"
Set-TextValue $ws.Range("E23") "  +0.63%  "

Set-TextValue $ws.Range("D24") "5.40"
Set-TextValue $ws.Range("E24") "  -12.57%  "

Set-TextValue $ws.Range("D25") "56.35"
Set-TextValue $ws.Range("E25") "  -9.97%  "

Set-TextValue $ws.Range("E26") "  +1.07%  "

Set-TextValue $ws.Range("D27") "0.386"
Set-TextValue $ws.Range("E27") "  -8.55%  "

Set-TextValue $ws.Range("D28") "2.552.37"
Set-TextValue $ws.Range("E28") "  -9.40%  "

Set-TextValue $ws.Range("E29") "  -8.37%  "

Set-TextValue $ws.Range("D30") "7.12"
Set-TextValue $ws.Range("E30") "  -3.37%  "

Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  +0.09%  "

Set-TextValue $ws.Range("E32") "  -12.92%  "

Set-TextValue $ws.Range("D33") "146.02"
Set-TextValue $ws.Range("E33") "  -2.63%  "

Set-TextValue $ws.Range("D34") "17.79"
Set-TextValue $ws.Range("E34") "  -6.38%  "

Set-TextValue $ws.Range("E35") "  -9.42%  "

Set-TextValue $ws.Range("E36") "  -5.93%  "

Set-TextValue $ws.Range("E37") "  -14.03%  "

Set-TextValue $ws.Range("E38") "  -5.12%  "

Set-TextValue $ws.Range("D39") "0.794"
Set-TextValue $ws.Range("E39") "  -14.37%  "

Set-TextValue $ws.Range("E40") "  +0.17%  "

Set-TextValue $ws.Range("D41") "32.82"
Set-TextValue $ws.Range("E41") "  -6.72%  "

Set-TextValue $ws.Range("D42") "0.597"
Set-TextValue $ws.Range("E42") "  +0.85%  "

Set-TextValue $ws.Range("D45") "10.12"
Set-TextValue $ws.Range("E45") "  -2.26%  "

Set-TextValue $ws.Range("D46") "1.24"
Set-TextValue $ws.Range("E46") "  -9.37%  "

Set-TextValue $ws.Range("D47") "1.935.20"
Set-TextValue $ws.Range("E47") "  -10.99%  "

Set-TextValue $ws.Range("D48") "0.0884"
Set-TextValue $ws.Range("E48") "  +0.07%  "

Set-TextValue $ws.Range("E49") "  -3.56%  "

Set-TextValue $ws.Range("D50") "232.58"
Set-TextValue $ws.Range("E50") "  +6.52%  "

Set-TextValue $ws.Range("D51") "16.59"
Set-TextValue $ws.Range("E51") "  -11.11%  "

# Row 43/44: swap Filecoin and Hedera, with updated values
Set-TextValue $ws.Range("B43") "Hedera"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D43") "0.0523"
Set-TextValue $ws.Range("E43") "  -5.32%  "

Set-TextValue $ws.Range("B44") "Filecoin"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "3.26"
Set-TextValue $ws.Range("E44") "  -7.92%  "
